$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column for "Fornecedor" right after "Nome do Cliente" (old col G) ---
$ws.Columns.Item(7).Insert()

# --- Insert a new column for "Quantidade" right after "Tipo de Movimentação" (old col W, now X) ---
$ws.Columns.Item(25).Insert()

# --- Populate the new header cells, matching the shared-string creation order seen in the target file ---
# (Valor ICMS, Valor PIS, Valor COFINS, Quantidade, Fornecedor)
$ws.Range("AB1").Value = "Valor ICMS"
$ws.Range("AC1").Value = "Valor PIS"
$ws.Range("AD1").Value = "Valor COFINS"
$ws.Range("Y1").Value = "Quantidade"
$ws.Range("G1").Value = "Fornecedor"

# --- Apply the same header style as the neighbouring "totals" block (bold + grey fill) to the 3 new trailing columns ---
$ws.Range("AA1").Copy()
$ws.Range("AB1:AD1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the values (PasteSpecial of formats only should not have touched them, but make sure)
$ws.Range("AB1").Value = "Valor ICMS"
$ws.Range("AC1").Value = "Valor PIS"
$ws.Range("AD1").Value = "Valor COFINS"

# --- Column widths for the newly introduced columns ---
$ws.Columns.Item(7).ColumnWidth = 29.833333333333332   # G  -> Fornecedor   (target xml width 30.7109375)
$ws.Columns.Item(25).ColumnWidth = 29.5                # Y  -> Quantidade   (target xml width 30.28515625)
$ws.Columns.Item(28).ColumnWidth = 9.833333333333334   # AB -> Valor ICMS   (target xml width 10.7109375)
$ws.Columns.Item(29).ColumnWidth = 8.0                 # AC -> Valor PIS    (target xml width 8.85546875)
$ws.Columns.Item(30).ColumnWidth = 11.833333333333334  # AD -> Valor COFINS (target xml width 12.7109375)

# --- Selection / view tidy-up to match the saved file ---
$ws.Range("G2").Select()
$ws.Application.ActiveWindow.ScrollColumn = 1
